# Edit script for "TEK-2000/Observer capabilities .docx"
# Applies the three visible-text changes described by the commit diff:
#   1. "Propeller size to use" -> "Propeller size to use (1206 as from kit)"
#   2. "100% at cruise/throttle high/throttle low " -> "100% rates at cruise/throttle high/throttle low "
#   3. "75% at cruise/throttle high/throttle low"  -> "75% rates at cruise/throttle high/throttle low"

$d = $word.ActiveDocument

# 1) Extend "Propeller size to use" with the kit note.
$d.Content.Find.Execute(
    "Propeller size to use",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Propeller size to use (1206 as from kit)",
    2) | Out-Null

# 2) Add "rates " before "at cruise/throttle high/throttle low " for the 100% bullet.
$d.Content.Find.Execute(
    "100% at cruise/throttle high/throttle low",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "100% rates at cruise/throttle high/throttle low",
    2) | Out-Null

# 3) Add "rates " before "at cruise/throttle high/throttle low" for the 75% bullet.
$d.Content.Find.Execute(
    "75% at cruise/throttle high/throttle low",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "75% rates at cruise/throttle high/throttle low",
    2) | Out-Null
